# Generate Report for Archive
# - Update the localization status text from "Ready for handoff" to "In Translation"
#   (appears on the Overview sheet, and on each per-language status sheet).
# - Narrow the per-language / status columns to match the refreshed report layout.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- 1. Replace the status text everywhere it is used ---
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count
    $rowStart = $used.Row
    $colStart = $used.Column

    for ($r = $rowStart; $r -lt ($rowStart + $rowCount); $r++) {
        for ($c = $colStart; $c -lt ($colStart + $colCount); $c++) {
            $cell = $ws.Cells.Item($r, $c)
            # NOTE: put the literal on the left -- if the cell holds a boolean,
            # "-eq" coerces the right-hand side to bool when the LHS is bool,
            # which would otherwise make every truthy cell match.
            if ($oldStatus -eq $cell.Value2) {
                $cell.Value = $newStatus
            }
        }
    }
}

# --- 2. Narrow the status columns ---
# Target layout width is ~13.41 characters; the engine's ColumnWidth setter
# only lands on a 1/6-character grid, so 12.5 is the input that snaps to the
# closest achievable stored width.
$newColumnWidth = 12.5

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth   # column E: zh-cn
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth   # column F: de-de

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth       # column C: Status

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth       # column C: Status
